$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update TPM-derived ligand/receptor expression and specificity metrics
# (columns G,H,I,J,M,N,O,P,Q,R,S,T) for rows 2-10, per new TPM values.
$ws.Range("G2").Value = 1.805635333333334
$ws.Range("H2").Value = 5.416906000000001
$ws.Range("I2").Value = 0.01900969238460649
$ws.Range("J2").Value = 0.01900969238460649
$ws.Range("M2").Value = 8.676671
$ws.Range("N2").Value = 26.030013
$ws.Range("O2").Value = 0.1325240072999665
$ws.Range("P2").Value = 0.1325240072999665
$ws.Range("Q2").Value = 15.66690373330867
$ws.Range("R2").Value = 141.002133599778
$ws.Range("S2").Value = 0.002519240612347708
$ws.Range("T2").Value = 0.002519240612347708
$ws.Range("G3").Value = 1.805635333333334
$ws.Range("H3").Value = 5.416906000000001
$ws.Range("I3").Value = 0.01900969238460649
$ws.Range("J3").Value = 0.01900969238460649
$ws.Range("M3").Value = 37.74750533333334
$ws.Range("O3").Value = 0.5765403197090441
$ws.Range("P3").Value = 0.576540319709044
$ws.Range("Q3").Value = 68.15822937505513
$ws.Range("R3").Value = 613.4240643754962
$ws.Range("S3").Value = 0.01095985412499161
$ws.Range("T3").Value = 0.0109598541249916
$ws.Range("G4").Value = 1.805635333333334
$ws.Range("H4").Value = 5.416906000000001
$ws.Range("I4").Value = 0.01900969238460649
$ws.Range("J4").Value = 0.01900969238460649
$ws.Range("M4").Value = 19.04827033333333
$ws.Range("N4").Value = 57.144811
$ws.Range("O4").Value = 0.2909356729909895
$ws.Range("P4").Value = 0.2909356729909895
$ws.Range("Q4").Value = 34.39422995275179
$ws.Range("R4").Value = 309.5480695747661
$ws.Range("S4").Value = 0.005530597647267177
$ws.Range("T4").Value = 0.005530597647267177
$ws.Range("I5").Value = 0.7995527524661065
$ws.Range("J5").Value = 0.7995527524661064
$ws.Range("M5").Value = 8.676671
$ws.Range("N5").Value = 26.030013
$ws.Range("O5").Value = 0.1325240072999665
$ws.Range("P5").Value = 0.1325240072999665
$ws.Range("Q5").Value = 658.9541665983024
$ws.Range("R5").Value = 5930.587499384721
$ws.Range("S5").Value = 0.1059599348045266
$ws.Range("T5").Value = 0.1059599348045266
$ws.Range("I6").Value = 0.7995527524661065
$ws.Range("J6").Value = 0.7995527524661064
$ws.Range("M6").Value = 37.74750533333334
$ws.Range("O6").Value = 0.5765403197090441
$ws.Range("P6").Value = 0.576540319709044
$ws.Range("Q6").Value = 2866.753380195198
$ws.Range("R6").Value = 25800.78042175678
$ws.Range("S6").Value = 0.4609743995310552
$ws.Range("T6").Value = 0.4609743995310551
$ws.Range("I7").Value = 0.7995527524661065
$ws.Range("J7").Value = 0.7995527524661064
$ws.Range("M7").Value = 19.04827033333333
$ws.Range("N7").Value = 57.144811
$ws.Range("O7").Value = 0.2909356729909895
$ws.Range("P7").Value = 0.2909356729909895
$ws.Range("Q7").Value = 1446.630522540366
$ws.Range("R7").Value = 13019.67470286329
$ws.Range("S7").Value = 0.2326184181305247
$ws.Range("T7").Value = 0.2326184181305247
$ws.Range("G8").Value = 17.23384333333334
$ws.Range("H8").Value = 51.70153000000001
$ws.Range("I8").Value = 0.1814375551492871
$ws.Range("J8").Value = 0.1814375551492871
$ws.Range("M8").Value = 8.676671
$ws.Range("N8").Value = 26.030013
$ws.Range("O8").Value = 0.1325240072999665
$ws.Range("P8").Value = 0.1325240072999665
$ws.Range("Q8").Value = 149.5323886688767
$ws.Range("R8").Value = 1345.79149801989
$ws.Range("S8").Value = 0.02404483188309219
$ws.Range("T8").Value = 0.02404483188309219
$ws.Range("G9").Value = 17.23384333333334
$ws.Range("H9").Value = 51.70153000000001
$ws.Range("I9").Value = 0.1814375551492871
$ws.Range("J9").Value = 0.1814375551492871
$ws.Range("M9").Value = 37.74750533333334
$ws.Range("O9").Value = 0.5765403197090441
$ws.Range("P9").Value = 0.576540319709044
$ws.Range("Q9").Value = 650.5345931388313
$ws.Range("R9").Value = 5854.811338249481
$ws.Range("S9").Value = 0.1046060660529973
$ws.Range("T9").Value = 0.1046060660529973
$ws.Range("G10").Value = 17.23384333333334
$ws.Range("H10").Value = 51.70153000000001
$ws.Range("I10").Value = 0.1814375551492871
$ws.Range("J10").Value = 0.1814375551492871
$ws.Range("M10").Value = 19.04827033333333
$ws.Range("N10").Value = 57.144811
$ws.Range("O10").Value = 0.2909356729909895
$ws.Range("P10").Value = 0.2909356729909895
$ws.Range("Q10").Value = 328.2749066956479
$ws.Range("R10").Value = 2954.474160260831
$ws.Range("S10").Value = 0.0527866572131976
$ws.Range("T10").Value = 0.0527866572131976
